# Add a new fee-earner mapping row (Zoe Baverstock) to the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "Zoe Baverstock"
$ws.Range("B17").Value = "Nicola.Daniel@taylorslegal.com"

# Match column B's existing "Hyperlink" cell style used by every other row.
$ws.Range("B17").Style = "Hyperlink"

# Leave the selection where the author's last save left it.
$ws.Range("E15").Select()
